# Update "想去人数" (attendance interest count) figures in column F
# for the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 7175
$ws1.Range("F4").Value  = 205
$ws1.Range("F5").Value  = 140
$ws1.Range("F6").Value  = 1098
$ws1.Range("F7").Value  = 177
$ws1.Range("F8").Value  = 7
$ws1.Range("F9").Value  = 65
$ws1.Range("F10").Value = 13

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 7175
$ws4.Range("F4").Value  = 205
$ws4.Range("F5").Value  = 140
$ws4.Range("F6").Value  = 1098
$ws4.Range("F7").Value  = 177
$ws4.Range("F9").Value  = 7
$ws4.Range("F10").Value = 65
$ws4.Range("F11").Value = 13
